# Documentacion de la primera iteracion
#
# The author bumped the "Factor de Productividad en hrs/cup" (C42 on the
# "Factor de complejidad Téc y Amb" sheet) from 4 to 8, which ripples
# through the "Estimación de Esfuerzo" sheet via its formulas, then left
# the workbook with the "Estimación de Esfuerzo" sheet active/selected.

$wb = $excel.ActiveWorkbook

# --- Core data edit -------------------------------------------------
$wsFactor = $wb.Worksheets.Item("Factor de complejidad Téc y Amb")
$wsFactor.Range("C42").Value = 8

# --- View / selection bookkeeping to match where the author ended up -
$wsInfGral = $wb.Worksheets.Item("Inf.Gral")
$wsInfGral.Activate()
$wsInfGral.Range("B5").Select()

$wsInstrucciones = $wb.Worksheets.Item("Instrucciones")
$wsInstrucciones.Activate()
$wsInstrucciones.Range("C50").Select()

$wsFactor.Activate()
$excel.ActiveWindow.Zoom = 85
$wsFactor.Range("C43").Select()

$wsTamano = $wb.Worksheets.Item("Estimación de Tamaño UCP")
$wsTamano.Activate()
$wsTamano.Range("F45").Select()

$wsRecursos = $wb.Worksheets.Item("Recursos")
$wsRecursos.Activate()
$wsRecursos.Range("C24").Select()

$wsEsfuerzo = $wb.Worksheets.Item("Estimación de Esfuerzo")
$wsEsfuerzo.Activate()
$wsEsfuerzo.Range("E31").Select()
